# json-example.xlsx: switch the "Cross Rate" sub-data lookup over to the
# "Main" endpoint with a KRW-USD key, and add a numeric-conversion row for
# the raw JSON string WEBSERVICE() returns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 is the header row; F5 used to read "Cross Rate" - rename the column
# to "Keys" since F6 below now holds a currency-pair key instead of a rate
# name.
$ws.Range("F5").Value = "Keys"

# Row 6 holds the actual parameter values used to build the request URL.
# F6 used to be the urlencoded "Cross%20Rate(US$)" sub-data selector; it's
# now the KRW-USD key used by the main-data endpoint.
$ws.Range("F6").Value = "KRW-USD"

# B8 builds the request URL. Switch from the Sub endpoint (D6) to the Main
# endpoint (C6) and drop the old "-1" suffix that was only needed for the
# sub-data call.
$ws.Range("B8").Formula = '=_xlfn.CONCAT(B6,C6,E6,F6)'

# B9 calls the web service with the rebuilt URL (formula unchanged, but it
# now resolves against the Main endpoint built above).
$ws.Range("B9").Formula = '=_xlfn.WEBSERVICE(B8)'

# New row 10: WEBSERVICE() comes back as a quoted JSON string (with a
# trailing newline), so strip the quotes/newline and coerce it to a real
# number for downstream use.
$ws.Range("B10").Formula = '=VALUE(SUBSTITUTE(SUBSTITUTE(B9, CHAR(34), ""), CHAR(10), ""))'

# Match the author's last selection before saving.
$ws.Range("B11").Select()
